$wb = $excel.ActiveWorkbook

$nodes = $wb.Worksheets.Item("nodes")
$interactions = $wb.Worksheets.Item("interactions")
$units = $wb.Worksheets.Item("units")

# --- nodes sheet: node 13 becomes a flier, node 18's z (D19) drops from 2 to 1 ---
$nodes.Range("G2").Copy($nodes.Range("G14"))
$nodes.Range("D19").Value = 1

# --- interactions sheet: fix bridge melee + LZ ("flier" vs "all") archer rules ---
$interactions.Range("C2").Value = 2
$interactions.Range("C3").Value = 2
$interactions.Range("C4").Value = 2

# give the bridge-to-bridge rows an (empty, styled) flier column cell
$interactions.Range("A1").Copy()
$interactions.Range("E5").PasteSpecial(-4122)
$interactions.Range("A1").Copy()
$interactions.Range("E6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# row 7 ("flier" -> "all") loses its flier fill and now reads all/all
$interactions.Range("A7").Value = "all"
$interactions.Range("A7").ClearFormats()

# drop the now-unused h3/h4 interaction rows (old rows 9-12); old row 13 (towers/other) becomes row 9
$interactions.Rows("9:12").Delete()

# --- restore cursor positions on each sheet, matching the saved view state ---
$nodes.Activate()
$nodes.Range("K27").Select()

$interactions.Activate()
$interactions.Range("E31").Select()

$units.Activate()
$units.Range("B41").Select()

$nodes.Activate()
